# Variables Funciones y Graficas.pptx - "cambios en los hiper vinculos"
#
# 1) Update the cached footer date text (02/02/2021 -> 03/02/2021) on the
#    slide master and every slide layout.
# 2) On slide 1 (title slide):
#    - shrink the title placeholder height
#    - resize/reposition the subtitle placeholder and drop its fixed
#      autofit scale back to a plain normAutofit
#    - change the last bullet to the "git clone" line, add a new
#      hyperlinked "Link de descarga del contenido" line (pointing at the
#      repo) and a trailing blank paragraph
#    - shrink/reposition the logo picture
#    - add a slow 2s slide transition that doesn't advance on click

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer date field text: 02/02/2021 -> 03/02/2021
# ---------------------------------------------------------------------
function Update-DateText($shape) {
    if ($shape.HasTextFrame -eq -1) {
        $tr = $shape.TextFrame.TextRange
        $found = $tr.Find("02/02/2021", 0)
        if ($found -ne $null) {
            $found.Text = "03/02/2021"
        }
    }
}

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateText($master.Shapes.Item($i))
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateText($layout.Shapes.Item($i))
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 edits
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# -- Title placeholder: only the height changes
$title = $s1.Shapes.Item(1)
$title.Height = 3090253 / 12700

# -- Subtitle placeholder: reposition/resize + simplify autofit
$sub = $s1.Shapes.Item(2)
$sub.Left = 685673 / 12700
$sub.Top = 3428999 / 12700
$sub.Width = 9228201 / 12700
$sub.Height = 2932043 / 12700
$sub.TextFrame.AutoSize = 2

$tr = $sub.TextFrame.TextRange
$lastPara = $tr.Paragraphs(4)
$lastPara.Text = "Git clone https://github.com/emena16/Topicos.git"
$tr.InsertAfter("`rLink de descarga del contenido`r")

$tr2 = $sub.TextFrame.TextRange
$linkPara = $tr2.Paragraphs(5)
$hlink = $linkPara.ActionSettings.Item(1).Hyperlink
$hlink.Address = "https://github.com/emena16/Topicos.git"

# -- Logo picture: shrink + reposition
$pic = $s1.Shapes.Item(3)
$pic.Left = 9713843 / 12700
$pic.Top = -4730 / 12700
$pic.Width = 2454149 / 12700
$pic.Height = 2454149 / 12700

# -- Slide transition: slow, 2s, no advance on click
$trans = $s1.SlideShowTransition
$trans.Duration = 2
$trans.Speed = [int][Microsoft.Office.Interop.PowerPoint.PpTransitionSpeed]::ppTransitionSpeedSlow
$trans.AdvanceOnClick = 0
